# "User Complaint filing backend updated"
# Column A holds user phone numbers used by the complaint-filing backend.
# Renumber A2:A23 to the new phone-number range. Rows 3,4,5,6,9,10,11,12,14
# previously carried "=<prev cell>+5" helper formulas; the new backend feed
# supplies literal values for every row, so those formulas are replaced with
# plain numbers (matching the look of the already-literal rows around them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New phone numbers for A2:A23 (literal values - no formulas anywhere in the column).
$phoneNumbers = @{
    2  = 9012345678
    3  = 9012345679
    4  = 9012345680
    5  = 9012345681
    6  = 9012345682
    7  = 9012345683
    8  = 9012345684
    9  = 9012345685
    10 = 9012345686
    11 = 9012345687
    12 = 9012345688
    13 = 9012345689
    14 = 9012345690
    15 = 9012345691
    16 = 9012345692
    17 = 9012345693
    18 = 9012345694
    19 = 9012345695
    20 = 9012345696
    21 = 9012345697
    22 = 9012345698
    23 = 9012345699
}

# Rows that used to hold a "=prevCell+5" formula - these need to become plain
# literal values styled like the rest of the column (Calibri 11 / black, the
# same look already used by A2, A7, A8, A13, A15:A23).
$formulaRows = 3,4,5,6,9,10,11,12,14

foreach ($r in 2..23) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $phoneNumbers[$r]
}

foreach ($r in $formulaRows) {
    $cell = $ws.Cells.Item($r, 1)
    # Match the literal-value style (font) already used elsewhere in the column.
    $cell.Font.Color = 0
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    # Row grows a touch taller to match the other literal rows once it carries
    # that font.
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Columns settled a little narrower once the new values/styles took effect.
$ws.Columns.Item(1).ColumnWidth = 13.748299319727868
$ws.Range("B1:AE1").EntireColumn.ColumnWidth = 7.534013605442177

# Selection left on the edited range.
$ws.Range("A2:A23").Select() | Out-Null
